$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new journal entry row (row 4): date + event label
$d = Get-Date -Year 2021 -Month 3 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Range("B4").Value = $d.Date
$ws.Range("C4").Value = "Finalisation de la grille"

# Update the table's visual style
$lo = $ws.ListObjects.Item(1)
$lo.TableStyle = "TableStyleLight2"

# Move the active selection to C4, matching the saved view state
$ws.Range("C4").Select()
